$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq "Administrator, Developer, Miss Dina Nasr") {
        $cell.Value = "Miss Dina Nasr, Administrator, Developer"
    }
    elseif ($val -eq "Administrator, Miss Dina Nasr") {
        $cell.Value = "Miss Dina Nasr, Administrator"
    }
}
